# Generate Report for handback
# Fills in the "Latest Target File" / "Latest Handback File" columns (E/F) for the
# two localized-file rows on the zh-cn and de-de sheets, flips the Status column
# text to reflect the handback, and (for de-de, which has actually been handed
# back) records the real "Latest Handback DateTime" in column G.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---- zh-cn sheet (Worksheets(2)) -----------------------------------------
$wsZh = $wb.Worksheets.Item(2)

$zhMd1  = "6dbb3ec4-4087-411a-9fa1-96ab82273acb.md"
$zhXlf1 = "6dbb3ec4-4087-411a-9fa1-96ab82273acb.e1431478200fed716e24059e35f65ce3e47d680e.zh-cn.xlf"
$zhMd2  = "dc09de62-9666-4302-a264-42801c2a4122.md"
$zhXlf2 = "dc09de62-9666-4302-a264-42801c2a4122.ab581e08832533e337b1eb5e1640b9344b36579b.zh-cn.xlf"

$wsZh.Range("B2").Value = $statusHandedBack
$wsZh.Range("B3").Value = $statusHandedBack

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/870965144f404e7eb8153f1f58263687d8130029/e2e/$zhMd1", "", "", $zhMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/85cc1b986092c529fe47758585d333f97f433c10/ol-handback/OpenLocalizationTest/oltest.zh-cn/xinjiang/ht/$zhXlf1", "", "", $zhXlf1)

$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/870965144f404e7eb8153f1f58263687d8130029/e2e/$zhMd2", "", "", $zhMd2)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/85cc1b986092c529fe47758585d333f97f433c10/ol-handback/OpenLocalizationTest/oltest.zh-cn/xinjiang/ht/$zhXlf2", "", "", $zhXlf2)

# ---- de-de sheet (Worksheets(3)) -----------------------------------------
$wsDe = $wb.Worksheets.Item(3)

$deMd1  = "6dbb3ec4-4087-411a-9fa1-96ab82273acb.md"
$deXlf1 = "6dbb3ec4-4087-411a-9fa1-96ab82273acb.e1431478200fed716e24059e35f65ce3e47d680e.de-de.xlf"
$deMd2  = "dc09de62-9666-4302-a264-42801c2a4122.md"
$deXlf2 = "dc09de62-9666-4302-a264-42801c2a4122.ab581e08832533e337b1eb5e1640b9344b36579b.de-de.xlf"

$wsDe.Range("B2").Value = $statusHandedBack
$wsDe.Range("B3").Value = $statusHandedBack

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/870965144f404e7eb8153f1f58263687d8130029/e2e/$deMd1", "", "", $deMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1d5da10e87466276effaa1a2d6627304df1702ab/ol-handback/OpenLocalizationTest/oltest.de-de/xinjiang/ht/$deXlf1", "", "", $deXlf1)

$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/870965144f404e7eb8153f1f58263687d8130029/e2e/$deMd2", "", "", $deMd2)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1d5da10e87466276effaa1a2d6627304df1702ab/ol-handback/OpenLocalizationTest/oltest.de-de/xinjiang/ht/$deXlf2", "", "", $deXlf2)

# de-de has actually completed its round trip back to en-US, so its "Latest
# Handback DateTime" column (G) now carries a real timestamp instead of the
# zero-date placeholder.
$wsDe.Range("G2").Value = "2016-01-26 03:43:05"
$wsDe.Range("G3").Value = "2016-01-26 03:43:05"
